$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Object Code")

# Update B6 to the new value "Item draw" (replaces "PowerUp draw")
$ws.Range("B6").Value = "Item draw"

# Delete row 7 entirely (A7=11, B7="PowerUp grabbed"), shifting cells up
$ws.Rows("7:7").Delete()
